{"js": "// Update the arithmetic problems (\"NN\u00f7N=\") in the first table of the\n// document. Each text-bearing row of the table (rows 0, 4, 8, 12, 16)\n// has 5 cells with a single \"a\u00f7b=\" expression; replace each with the\n// new expression from the commit, targeting cells by (row, col) so\n// that duplicate values (e.g. \"51\u00f79=\" appearing twice) are replaced\n// independently and correctly.\n\nconst replacements = [\n  { row: 0, col: 0, oldText: \"41\u00f74=\", newText: \"34\u00f73=\" },\n  { row: 0, col: 1, oldText: \"45\u00f76=\", newText: \"75\u00f72=\" },\n  { row: 0, col: 2, oldText: \"74\u00f79=\", newText: \"13\u00f72=\" },\n  { row: 0, col: 3, oldText: \"47\u00f79=\", newText: \"90\u00f76=\" },\n  { row: 0, col: 4, oldText: \"69\u00f74=\", newText: \"14\u00f74=\" },\n\n  { row: 4, col: 0, oldText: \"27\u00f76=\", newText: \"75\u00f77=\" },\n  { row: 4, col: 1, oldText: \"74\u00f78=\", newText: \"99\u00f73=\" },\n  { row: 4, col: 2, oldText: \"77\u00f76=\", newText: \"85\u00f77=\" },\n  { row: 4, col: 3, oldText: \"64\u00f73=\", newText: \"86\u00f73=\" },\n  { row: 4, col: 4, oldText: \"56\u00f74=\", newText: \"84\u00f75=\" },\n\n  { row: 8, col: 0, oldText: \"42\u00f72=\", newText: \"92\u00f75=\" },\n  { row: 8, col: 1, oldText: \"98\u00f77=\", newText: \"58\u00f78=\" },\n  { row: 8, col: 2, oldText: \"75\u00f79=\", newText: \"85\u00f77=\" },\n  { row: 8, col: 3, oldText: \"26\u00f75=\", newText: \"45\u00f79=\" },\n  { row: 8, col: 4, oldText: \"51\u00f79=\", newText: \"15\u00f72=\" },\n\n  { row: 12, col: 0, oldText: \"70\u00f76=\", newText: \"28\u00f74=\" },\n  { row: 12, col: 1, oldText: \"51\u00f79=\", newText: \"68\u00f79=\" },\n  { row: 12, col: 2, oldText: \"29\u00f79=\", newText: \"20\u00f72=\" },\n  { row: 12, col: 3, oldText: \"73\u00f77=\", newText: \"46\u00f77=\" },\n  { row: 12, col: 4, oldText: \"34\u00f72=\", newText: \"10\u00f73=\" },\n\n  { row: 16, col: 0, oldText: \"13\u00f74=\", newText: \"57\u00f76=\" },\n  { row: 16, col: 1, oldText: \"64\u00f74=\", newText: \"45\u00f78=\" },\n  { row: 16, col: 2, oldText: \"29\u00f78=\", newText: \"99\u00f77=\" },\n  { row: 16, col: 3, oldText: \"93\u00f76=\", newText: \"51\u00f76=\" },\n  { row: 16, col: 4, oldText: \"50\u00f73=\", newText: \"16\u00f73=\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Collect the search result ranges for every cell first.\nconst searchResults = replacements.map((r) => {\n  const cell = table.getCell(r.row, r.col);\n  return cell.body.search(r.oldText, { matchCase: true, matchWholeWord: true });\n});\nsearchResults.forEach((res) => res.load(\"items\"));\nawait context.sync();\n\n// Replace the (single, exact) match found in each targeted cell.\nsearchResults.forEach((res, i) => {\n  const range = res.items[0];\n  range.insertText(replacements[i].newText, Word.InsertLocation.replace);\n});\nawait context.sync();\n", "ps1": "# Update the arithmetic problems (\"NN\u00f7N=\") in the first table of the\n# document. Each text-bearing row of the table (Word COM rows are\n# 1-indexed: 1, 5, 9, 13, 17) has 5 cells with a single \"a\u00f7b=\"\n# expression; replace each with the new expression from the commit,\n# targeting cells by (row, col) so that duplicate values (e.g.\n# \"51\u00f79=\" appearing twice) are replaced independently and correctly.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1;  Col = 1; Old = \"41\u00f74=\"; New = \"34\u00f73=\" },\n    @{ Row = 1;  Col = 2; Old = \"45\u00f76=\"; New = \"75\u00f72=\" },\n    @{ Row = 1;  Col = 3; Old = \"74\u00f79=\"; New = \"13\u00f72=\" },\n    @{ Row = 1;  Col = 4; Old = \"47\u00f79=\"; New = \"90\u00f76=\" },\n    @{ Row = 1;  Col = 5; Old = \"69\u00f74=\"; New = \"14\u00f74=\" },\n\n    @{ Row = 5;  Col = 1; Old = \"27\u00f76=\"; New = \"75\u00f77=\" },\n    @{ Row = 5;  Col = 2; Old = \"74\u00f78=\"; New = \"99\u00f73=\" },\n    @{ Row = 5;  Col = 3; Old = \"77\u00f76=\"; New = \"85\u00f77=\" },\n    @{ Row = 5;  Col = 4; Old = \"64\u00f73=\"; New = \"86\u00f73=\" },\n    @{ Row = 5;  Col = 5; Old = \"56\u00f74=\"; New = \"84\u00f75=\" },\n\n    @{ Row = 9;  Col = 1; Old = \"42\u00f72=\"; New = \"92\u00f75=\" },\n    @{ Row = 9;  Col = 2; Old = \"98\u00f77=\"; New = \"58\u00f78=\" },\n    @{ Row = 9;  Col = 3; Old = \"75\u00f79=\"; New = \"85\u00f77=\" },\n    @{ Row = 9;  Col = 4; Old = \"26\u00f75=\"; New = \"45\u00f79=\" },\n    @{ Row = 9;  Col = 5; Old = \"51\u00f79=\"; New = \"15\u00f72=\" },\n\n    @{ Row = 13; Col = 1; Old = \"70\u00f76=\"; New = \"28\u00f74=\" },\n    @{ Row = 13; Col = 2; Old = \"51\u00f79=\"; New = \"68\u00f79=\" },\n    @{ Row = 13; Col = 3; Old = \"29\u00f79=\"; New = \"20\u00f72=\" },\n    @{ Row = 13; Col = 4; Old = \"73\u00f77=\"; New = \"46\u00f77=\" },\n    @{ Row = 13; Col = 5; Old = \"34\u00f72=\"; New = \"10\u00f73=\" },\n\n    @{ Row = 17; Col = 1; Old = \"13\u00f74=\"; New = \"57\u00f76=\" },\n    @{ Row = 17; Col = 2; Old = \"64\u00f74=\"; New = \"45\u00f78=\" },\n    @{ Row = 17; Col = 3; Old = \"29\u00f78=\"; New = \"99\u00f77=\" },\n    @{ Row = 17; Col = 4; Old = \"93\u00f76=\"; New = \"51\u00f76=\" },\n    @{ Row = 17; Col = 5; Old = \"50\u00f73=\"; New = \"16\u00f73=\" }\n)\n\nforeach ($r in $replacements) {\n    $cellRange = $tbl.Cell($r.Row, $r.Col).Range\n    $cellRange.Find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 1)\n}\n"}
